$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2025-04-13 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-14 Monday", 2)

# The worksheet answers live in a single table; rows with content are
# 1-based rows 1, 5, 9, 13, 17 (0-based 0, 4, 8, 12, 16), 5 columns each.
$t = $d.Tables.Item(1)

function Replace-CellText($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # wdReplaceOne (1) -- some of the answer strings repeat verbatim in other
    # cells, so a scoped "replace all" (2) must not be used here or it will
    # replace every matching occurrence in the whole document instead of
    # just the text inside this cell.
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
}

# Row 1 (0-based row 0)
Replace-CellText $t 1 1 "49÷8=6, 1" "47÷4=11, 3"
Replace-CellText $t 1 2 "65÷9=7, 2" "31÷2=15, 1"
Replace-CellText $t 1 3 "78÷4=19, 2" "42÷4=10, 2"
Replace-CellText $t 1 4 "55÷4=13, 3" "94÷9=10, 4"
Replace-CellText $t 1 5 "96÷2=48, 0" "36÷7=5, 1"

# Row 5 (0-based row 4)
Replace-CellText $t 5 1 "14÷3=4, 2" "70÷8=8, 6"
Replace-CellText $t 5 2 "68÷8=8, 4" "61÷2=30, 1"
Replace-CellText $t 5 3 "29÷3=9, 2" "85÷3=28, 1"
Replace-CellText $t 5 4 "58÷5=11, 3" "31÷5=6, 1"
Replace-CellText $t 5 5 "34÷8=4, 2" "20÷2=10, 0"

# Row 9 (0-based row 8)
Replace-CellText $t 9 1 "53÷2=26, 1" "56÷2=28, 0"
Replace-CellText $t 9 2 "53÷4=13, 1" "25÷6=4, 1"
Replace-CellText $t 9 3 "83÷4=20, 3" "69÷7=9, 6"
Replace-CellText $t 9 4 "81÷9=9, 0" "47÷2=23, 1"
Replace-CellText $t 9 5 "34÷7=4, 6" "24÷3=8, 0"

# Row 13 (0-based row 12)
Replace-CellText $t 13 1 "29÷7=4, 1" "84÷3=28, 0"
Replace-CellText $t 13 2 "67÷7=9, 4" "63÷2=31, 1"
Replace-CellText $t 13 3 "53÷8=6, 5" "57÷4=14, 1"
Replace-CellText $t 13 4 "95÷7=13, 4" "79÷8=9, 7"
Replace-CellText $t 13 5 "40÷8=5, 0" "94÷9=10, 4"

# Row 17 (0-based row 16)
Replace-CellText $t 17 1 "28÷8=3, 4" "34÷9=3, 7"
Replace-CellText $t 17 2 "29÷7=4, 1" "15÷9=1, 6"
Replace-CellText $t 17 3 "14÷2=7, 0" "50÷8=6, 2"
Replace-CellText $t 17 4 "11÷2=5, 1" "53÷8=6, 5"
Replace-CellText $t 17 5 "70÷7=10, 0" "93÷8=11, 5"

Write-Output "Done"
